# Add the new "Reserved" (column E) marker "Y" to the last five product
# rows of the wishlist sheet, then leave the selection on the last cell
# touched (E20), matching the author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E20").Value = "Y"
$ws.Range("E21").Value = "Y"
$ws.Range("E22").Value = "Y"
$ws.Range("E23").Value = "Y"
$ws.Range("E24").Value = "Y"

$ws.Range("E20").Select()
